# Added send-approval serial & parallel cycle designer test-script data:
# a new "Valid User PA Name" row (row 9) on the Admin sheet, holding a PA
# user name, email and a hyperlinked password cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 9: Valid User PA Name / testuserpa@gmail.com / Test@123 -----------
$ws.Range("A9").Value = "Valid User PA Name"

$ws.Range("B9").Value = "testuserpa@gmail.com"
# New small "Prime" font for the PA user-name cell. Start from the closest
# existing style (C8: Arial 12 FF222222) so only Size/Name need changing,
# keeping the number of newly minted font records to a minimum.
$ws.Range("C8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Font.Size = 8
$ws.Range("B9").Font.Name = "Prime"

$ws.Range("C9").Value = "Test@123"
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:testuserpa@gmail.com")

$ws.Range("C9").Select()
